$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (A36) repeating the same keyword as A35 ("yoga outfit set"),
# matching the existing shared-string entry already used by A35.
$ws.Range("A36").Value = "yoga outfit set"

# Move/update the active selection to the newly added cell, as recorded in
# the saved sheet view (activeCell="A36" sqref="A36").
$ws.Range("A36").Select()
